{"js": "// Remove extra (empty) table rows from the header/contact-info table.\n//\n// The document has a single table whose rows hold contact-info lines\n// (email, city, social links). A few rows in the middle are blank\n// placeholder rows (a \"Compact\" styled paragraph with no text and no\n// hyperlinks) that need to be deleted, while rows that carry real\n// content - including the trailing row that only contains hyperlinks\n// separated by single-space runs - must be preserved.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  return \"no tables found\";\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather each row's text and whether it contains any hyperlinks.\nconst rowRanges = rows.items.map((row) => row.getRange());\nrowRanges.forEach((r) => r.load(\"text\"));\nconst hyperlinkCollections = rowRanges.map((r) => r.getHyperlinkRanges());\nhyperlinkCollections.forEach((h) => h.load(\"items\"));\nawait context.sync();\n\n// Delete from the bottom up so earlier indices stay valid.\nfor (let i = rows.items.length - 1; i >= 0; i--) {\n  const isBlankText = rowRanges[i].text.trim().length === 0;\n  const hasHyperlinks = hyperlinkCollections[i].items.length > 0;\n\n  if (isBlankText && !hasHyperlinks) {\n    rows.items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove extra (empty) table rows from the header/contact-info table.\n#\n# The document has a single table whose rows hold contact-info lines\n# (email, city, social links). A few rows in the middle are blank\n# placeholder rows (a \"Compact\" styled paragraph with no text and no\n# hyperlinks) that need to be deleted, while rows that carry real\n# content - including the trailing row that only contains hyperlinks\n# separated by single-space runs - must be preserved.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Word's Range.Text includes trailing end-of-cell/end-of-row markers\n# (CR + BEL, char codes 13 and 7) even for \"empty\" cells, so strip those\n# along with plain whitespace before checking for real content.\nfor ($i = $table.Rows.Count; $i -ge 1; $i--) {\n    $row = $table.Rows.Item($i)\n    $rowRange = $row.Range\n    $trimmedText = $rowRange.Text.Trim([char]13, [char]7, [char]32, [char]9)\n    $isBlankText = ($trimmedText.Length -eq 0)\n    $hasHyperlinks = ($rowRange.Hyperlinks.Count -gt 0)\n\n    if ($isBlankText -and -not $hasHyperlinks) {\n        $row.Delete()\n    }\n}\n"}
